$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell 'D2' '66.232.74'
Set-TextCell 'E2' '  -0.53%  '

Set-TextCell 'D3' '3.535.70'
Set-TextCell 'E3' '  -1.74%  '

Set-TextCell 'D4' '0.998'
Set-TextCell 'E4' '  -0.18%  '

Set-TextCell 'D5' '603.54'
Set-TextCell 'E5' '  -0.94%  '

Set-TextCell 'D6' '144.22'
Set-TextCell 'E6' '  -2.79%  '

Set-TextCell 'D7' '3.541.50'
Set-TextCell 'E7' '  -1.54%  '

Set-TextCell 'D8' '0.998'
Set-TextCell 'E8' '  -0.23%  '

Set-TextCell 'D10' '7.86'
Set-TextCell 'E10' '  -1.96%  '

Set-TextCell 'D11' '0.131'
Set-TextCell 'E11' '  -3.92%  '

Set-TextCell 'D12' '0.409'
Set-TextCell 'E12' '  -1.45%  '

Set-TextCell 'D13' '4.132.67'
Set-TextCell 'E13' '  -1.83%  '

Set-TextCell 'D14' '0.0000195'
Set-TextCell 'E14' '  -7.17%  '

Set-TextCell 'D15' '28.52'
Set-TextCell 'E15' '  -4.97%  '

Set-TextCell 'D16' '3.528.76'
Set-TextCell 'E16' '  -2.14%  '

Set-TextCell 'E17' '  +0.94%  '

Set-TextCell 'D18' '65.929.43'
Set-TextCell 'E18' '  -1.12%  '

Set-TextCell 'D19' '11.07'
Set-TextCell 'E19' '  -3.87%  '

Set-TextCell 'D20' '6.17'
Set-TextCell 'E20' '  -2.89%  '

Set-TextCell 'D21' '14.63'
Set-TextCell 'E21' '  -2.88%  '

Set-TextCell 'D22' '423.17'
Set-TextCell 'E22' '  -1.31%  '

Set-TextCell 'D23' '0.595'
Set-TextCell 'E23' '  -4.12%  '

Set-TextCell 'D24' '77.06'
Set-TextCell 'E24' '  -2.57%  '

Set-TextCell 'D25' '3.671.17'
Set-TextCell 'E25' '  -2.00%  '

Set-TextCell 'E26' '  +0.06%  '

Set-TextCell 'D27' '0.0000115'
Set-TextCell 'E27' '  -6.16%  '

Set-TextCell 'E28' '  -2.34%  '

Set-TextCell 'D29' '7.84'
Set-TextCell 'E29' '  -5.41%  '

Set-TextCell 'D30' '8.90'
Set-TextCell 'E30' '  -4.44%  '

Set-TextCell 'D31' '0.998'
Set-TextCell 'E31' '  -0.20%  '

Set-TextCell 'D32' '3.538.38'
Set-TextCell 'E32' '  -1.63%  '

Set-TextCell 'D33' '0.155'
Set-TextCell 'E33' '  -1.03%  '

Set-TextCell 'D34' '24.26'
Set-TextCell 'E34' '  -4.97%  '

Set-TextCell 'E35' '  -0.04%  '

Set-TextCell 'D36' '1.35'
Set-TextCell 'E36' '  -7.42%  '

Set-TextCell 'D37' '7.61'
Set-TextCell 'E37' '  -3.22%  '

Set-TextCell 'B38' 'ImmutableX'
Set-TextCell 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D38' '1.64'
Set-TextCell 'E38' '  -4.58%  '

Set-TextCell 'B39' 'Monero'
Set-TextCell 'C39' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D39' '176.58'
Set-TextCell 'E39' '  -0.20%  '

Set-TextCell 'D40' '5.24'
Set-TextCell 'E40' '  -7.19%  '

Set-TextCell 'D41' '0.0820'
Set-TextCell 'E41' '  -4.59%  '

Set-TextCell 'D42' '0.861'
Set-TextCell 'E42' '  -4.22%  '

Set-TextCell 'D43' '4.98'
Set-TextCell 'E43' '  -4.97%  '

Set-TextCell 'D44' '45.40'
Set-TextCell 'E44' '  -1.79%  '

Set-TextCell 'D45' '1.77'
Set-TextCell 'E45' '  -7.93%  '

Set-TextCell 'D46' '0.999'
Set-TextCell 'E46' '  -0.13%  '

Set-TextCell 'D47' '2.40'
Set-TextCell 'E47' '  -7.60%  '

Set-TextCell 'D48' '23.91'
Set-TextCell 'E48' '  -1.74%  '

Set-TextCell 'D49' '7.06'
Set-TextCell 'E49' '  -1.93%  '

Set-TextCell 'D50' '1.13'
Set-TextCell 'E50' '  -5.82%  '

Set-TextCell 'D51' '0.910'
Set-TextCell 'E51' '  -4.42%  '
